$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 6 date
$ws.Range("D6").Value = 44642

# Update existing row 7: date and volume
$ws.Range("D7").Value = 44637
$ws.Range("J7").Value = 100

# Insert new row 8, replicating the original (pre-edit) row 7 contents
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 44628
$ws.Range("D8").NumberFormat = $ws.Range("D7").NumberFormat
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100112043
$ws.Range("G8").Value = "Pepino dulce"
$ws.Range("H8").Value = "Cultivar IV Región"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15500
$ws.Range("N8").Value = "$/bandeja 18 kilos"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 861
$ws.Range("Q8").Value = 18
$ws.Range("R8").Value = "Hortaliza"
